# Update the "dSF" column (F) values for the affected rows,
# per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 6
$ws.Range("F3").Value  = -2
$ws.Range("F4").Value  = 7
$ws.Range("F5").Value  = 2
$ws.Range("F6").Value  = -1
$ws.Range("F8").Value  = -2
$ws.Range("F9").Value  = 3
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = 3
